$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the McKinnon row, the duplicate Mordialloc (23/12) row, and the Southbank row,
# working from the bottom up so row numbers of earlier rows stay valid.
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(12).Delete()

# Update the Hampton exposure period.
$ws.Range("C10").Value = "28/12/20 12:50pm-2:40pm"

# Update the (now shifted-up) Melbourne / Left Bank row with its new site name and note.
$ws.Range("B12").Value = "Left Bank Melbourne Restaurant and Cocktail Bar, 1 Southbank Blvd"
$ws.Range("D12").Value = "Case attended bar"
